$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Hunk 1: "pass varchar(50) NOT NULL," -> "pass" + "word" + " varchar(" +
#         "255" + ") NOT NULL," (five separate runs), i.e. the USERS
#         table's password column is renamed/widened.
# ------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("pass varchar(50) NOT NULL,", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    [void]$rng1.Expand(4)   # wdParagraph - grab the whole paragraph (incl. leading/trailing tab runs)

    $xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="26B0077E" w14:textId="3DA09D4A" w:rsidR="002A4A35" w:rsidRDefault="002D1C78" w:rsidP="001631BF"><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:tab/><w:t>pass</w:t></w:r><w:r><w:t>word</w:t></w:r><w:r><w:t xml:space="preserve"> varchar(</w:t></w:r><w:r><w:t>255</w:t></w:r><w:r><w:t>) NOT NULL,</w:t></w:r><w:r w:rsidR="002A4A35"><w:tab/></w:r></w:p>'
    [void]$rng1.InsertXML($xml1)
}

# ------------------------------------------------------------------
# Hunk 2: the ISSUES table's trailing "version " / "serial" runs get
#         merged into a single "version serial" run. This is the last
#         "version" occurrence in the document (it sits right after the
#         uniquely-named ISSUES table), so anchor the search there.
# ------------------------------------------------------------------
$anchor = $d.Content
$foundAnchor = $anchor.Find.Execute("ISSUES", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundAnchor) {
    $rng2 = $d.Range($anchor.End, $d.Content.End)
    $found2 = $rng2.Find.Execute("version serial", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found2) {
        # Force the run merge: Word (and this host) only coalesces runs
        # when the text inside the range actually changes, so flip it to
        # a placeholder and back to the desired final text.
        $rng2.Text = "version serial__TMP__"
        $rng2b = $d.Content
        [void]$rng2b.Find.Execute("version serial__TMP__", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
        $rng2b.Text = "version serial"
    }
}
